# "Dissertation Script cleaning" pass on the pilot egg-size data.
#
# The egg_volume column (E) had been computed with the formula for the
# volume of a SPHERE: (4/3)*pi*r^3, using the raw length/width
# measurements (C = length, D = width) as if they were diameters, i.e.
# (PI()*4/3)*(C*D*D). Eggs are better approximated as a prolate spheroid,
# whose volume is (pi/6)*length*width^2, so every formula in column E is
# corrected to (PI()*1/6)*(C*D*D), keeping the same C/D cell references
# per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 361

# Row 2 carries its own (non-shared) formula, matching the workbook's
# existing convention of keeping the header-adjacent row separate from
# the big shared-formula block that follows it.
$ws.Range("E$firstDataRow").Formula = "=(PI()* 1/6)* (C$firstDataRow*D$firstDataRow*D$firstDataRow)"

# The remaining rows (3-361) share the corrected formula, relative
# references adjusting per row just like the original sheet.
$sharedFirstRow = $firstDataRow + 1
$ws.Range("E$sharedFirstRow`:E$lastDataRow").Formula = "=(PI()* 1/6)* (C$sharedFirstRow*D$sharedFirstRow*D$sharedFirstRow)"

# Reflect where the editor ended up after the cleaning pass: scrolled to
# the bottom of the data with the second-to-last egg_volume cell selected.
$ws.Range("E360").Select()
